# Handback status report regeneration: refresh the handoff/handback
# timestamps and the zh-cn/de-de "Priority" values (ht -> mt) that a
# fresh CI run of the Generate Report step would produce.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet --------------------------------------------------
# "Latest HO Xliff Generate Date" for the 723195df*.md row (rows 2 & 4,
# zh-cn + de-de columns share the same handoff timestamp string).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-20 10:15:18"
$wsOverview.Range("G4").Value = "2016-08-20 10:15:18"

# --- "zh-cn" sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority: ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-20 10:15:14"
$wsZhCn.Range("H4").Value = "2016-08-20 10:15:14"
# Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-20 10:15:31"
$wsZhCn.Range("K4").Value = "2016-08-20 10:15:31"

# --- "de-de" sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority: ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime (same value as the Overview sheet's
# handoff timestamp)
$wsDeDe.Range("H2").Value = "2016-08-20 10:15:18"
$wsDeDe.Range("H4").Value = "2016-08-20 10:15:18"
# Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-20 10:15:37"
$wsDeDe.Range("K4").Value = "2016-08-20 10:15:37"
